# Applies updated triangulation2 script values to the Mic Pair actual/estimated columns (B:G) for rows 2-11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6603181522184486
$ws.Range("C2").Value = 1.065759637188209
$ws.Range("D2").Value = 0.9112048919340322
$ws.Range("E2").Value = 0.6575963718820862
$ws.Range("F2").Value = 1.071598303772011
$ws.Range("G2").Value = 0.9070294784580499
$ws.Range("B3").Value = 0.3950008508465115
$ws.Range("C3").Value = 0.1587301587301587
$ws.Range("D3").Value = -0.4402577480568004
$ws.Range("E3").Value = 0.4081632653061225
$ws.Range("F3").Value = 0.1559849851160945
$ws.Range("G3").Value = -0.4308390022675737
$ws.Range("B4").Value = -0.4767228451427073
$ws.Range("C4").Value = -1.26984126984127
$ws.Range("D4").Value = -2.031350183375844
$ws.Range("E4").Value = -0.4761904761904762
$ws.Range("F4").Value = -1.26818552774021
$ws.Range("G4").Value = -2.040816326530612
$ws.Range("B5").Value = 0.6237795558828675
$ws.Range("C5").Value = 0.5895691609977325
$ws.Range("D5").Value = -0.002692718193406324
$ws.Range("E5").Value = 0.6122448979591837
$ws.Range("F5").Value = 0.6050531937043389
$ws.Range("G5").Value = 0
$ws.Range("B6").Value = 0.7086389277650147
$ws.Range("C6").Value = 1.292517006802721
$ws.Range("D6").Value = 1.541215953830037
$ws.Range("E6").Value = 0.7029478458049886
$ws.Range("F6").Value = 1.304139074833224
$ws.Range("G6").Value = 1.541950113378685
$ws.Range("B7").Value = 0.2473850426869561
$ws.Range("C7").Value = 0.06802721088435373
$ws.Range("D7").Value = -0.3802581138671877
$ws.Range("E7").Value = 0.2494331065759637
$ws.Range("F7").Value = 0.08084678929578659
$ws.Range("G7").Value = -0.3854875283446712
$ws.Range("B8").Value = 0.659905908185257
$ws.Range("C8").Value = 0.9297052154195011
$ws.Range("D8").Value = 0.5145133491967829
$ws.Range("E8").Value = 0.6575963718820862
$ws.Range("F8").Value = 0.9324018116487143
$ws.Range("G8").Value = 0.5215419501133787
$ws.Range("B9").Value = 0.5137392823977042
$ws.Range("C9").Value = 0.7256235827664399
$ws.Range("D9").Value = 0.4907005132262501
$ws.Range("E9").Value = 0.5215419501133787
$ws.Range("F9").Value = 0.7141090169048173
$ws.Range("G9").Value = 0.4988662131519275
$ws.Range("B10").Value = 0.6261701570170065
$ws.Range("C10").Value = 1.08843537414966
$ws.Range("D10").Value = 1.165374245474661
$ws.Range("E10").Value = 0.6349206349206349
$ws.Range("F10").Value = 1.07619612623426
$ws.Range("G10").Value = 1.179138321995465
$ws.Range("B11").Value = 0.7478644345922981
$ws.Range("C11").Value = 1.26984126984127
$ws.Range("D11").Value = 0.8777519492556363
$ws.Range("E11").Value = 0.7482993197278912
$ws.Range("F11").Value = 1.273801470238159
$ws.Range("G11").Value = 0.8843537414965986
